# The "Förändrad" (Changed) date column (C) for every data row on the
# "Avverkningsanmälningar" sheet is bumped by one day: serial date 45188
# (2023-09-19) becomes 45189 (2023-09-20). The data rows run from row 2
# through row 359 (dimension A1:Y359, header in row 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C359").Value = 45189
